$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Weekly driver report update for 2025-04-21
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 1219
$ws.Range("B4").Value = 41
$ws.Range("C4").Value = 1219
